$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(2).Collapsed = $true
